$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# The sheet is protected; lift protection so rows/cells can be edited, then
# restore the exact same protection settings at the end.
# ---------------------------------------------------------------------------
$ws.Unprotect("EF56")

# ---------------------------------------------------------------------------
# 1. Fill the previously-empty "Start" column (K) for the existing 16 task
#    rows (6-21) with 0.
# ---------------------------------------------------------------------------
for ($r = 6; $r -le 21; $r++) {
    $ws.Cells.Item($r, 11).Value = 0   # column K
}

# ---------------------------------------------------------------------------
# 2. Insert 10 new rows right before row 22 (the "Ideal burndown" summary
#    row), pushing everything from row 22 down to row 32.
# ---------------------------------------------------------------------------
$ws.Rows("22:31").Insert()

# Give the freshly-inserted rows the same formatting (fills/borders/etc.)
# used by the existing task rows (B:D = task-label style, E:Q = value style).
$ws.Range("B21:Q21").Copy() | Out-Null
$ws.Range("B22:Q31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Populate the new rows with the "Re-Checking / Re-Correct / Re-Designing"
#    logbook entries.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "Android"
$ws.Range("C22").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D22").Value = "Meringkas/menghapus topik Pengenalan Android "
$ws.Range("E22:K22").Value = 3

$ws.Range("B23").Value = "Android"
$ws.Range("C23").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D23").Value = "Meringkas/menghapus topik Trik Android"

$ws.Range("B24").Value = "Android"
$ws.Range("C24").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D24").Value = "Meringkas/menghapus topik  Masalah di Android"
$ws.Range("E24:K24").Value = 2

$ws.Range("B25").Value = "Android"
$ws.Range("C25").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D25").Value = "Meringkas/menghapus topik Benchmark"

$ws.Range("B26").Value = "Android"
$ws.Range("C26").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D26").Value = "Meringkas/menghapus topik Aplikasi Android yang Disarankan"

$ws.Range("B27").Value = "Android"
$ws.Range("C27").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D27").Value = "Meringkas/menghapus topik Android Google Play"
$ws.Range("E27:K27").Value = 2

$ws.Range("B28").Value = "Android"
$ws.Range("B29").Value = "Android"
$ws.Range("B30").Value = "Android"
$ws.Range("B31").Value = "Android"

# ---------------------------------------------------------------------------
# 4. The "Ideal burndown" (now row 32) and "Actual burndown" (now row 33)
#    totals need to widen their SUM range to include the new task row 22.
# ---------------------------------------------------------------------------
$ws.Range("E32").Formula = "=SUM(E6:E22)"
$ws.Range("E33").Formula = "=SUM(E6:E22)"

# ---------------------------------------------------------------------------
# 5. Update the chart's series source ranges to the new summary rows and
#    move the chart down by the height of the 10 inserted rows.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = '=SERIES("Ideal burndown",Sheet1!$F$5:$Q$5,Sheet1!$F$32:$Q$32,1)'
$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = '=SERIES("Actual burndown",Sheet1!$F$5:$Q$5,Sheet1!$F$33:$Q$33,2)'

$rowHeight = 15
$co.Top = $co.Top + ($rowHeight * 10)

# ---------------------------------------------------------------------------
# 6. Misc cosmetic sheet updates captured by the diff: selection/scroll
#    position and a couple of column widths.
# ---------------------------------------------------------------------------
$ws.Range("B13").Select()
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E27:K27").Select()

$ws.Columns("B:B").ColumnWidth = 18.42578125
$ws.Columns("C:C").ColumnWidth = 22.85546875

# ---------------------------------------------------------------------------
# 7. Re-apply the original sheet protection.
# ---------------------------------------------------------------------------
$ws.Protect("EF56", $false, $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true, $false, $false)
